# Generate Report for Handoff
# Update status + handoff timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest Handoff Date / Datetime timestamps bumped forward (stored as text)
$wsOverview.Range("D2").Value = "2016-28-19 02:28:51"
$wsZhCn.Range("E2").Value = "2016-03-19 02:28:48"
$wsDeDe.Range("E2").Value = "2016-03-19 02:28:51"
